$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B11 (Shelving Rack) to a formula 4*4 = 16
$ws.Range("B11").Formula = "=4*4"

# Update B14 (Clip Strip) to a formula 2*4 = 8
$ws.Range("B14").Formula = "=2*4"

# Update selection to C8
$ws.Range("C8").Select()
